{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// The instruction document told students to load data \"fr\u00e5n mappen\n// AdventureWorksDW\" (from the folder AdventureWorksDW). That folder\n// name was wrong; the commit corrects the wording to say the data\n// should be loaded from the \"data\" folder, which is a sub-folder of\n// \"Deluppgift 1\".\n//\n// Old: \"... fr\u00e5n mappen AdventureWorksDW. Transformera data ...\"\n// New: \"... fr\u00e5n mappen data som \u00e4r en undermapp under mappen\n//       Deluppgift 1. Transformera data ...\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldFragment = \"mappen AdventureWorksDW.\";\nconst newFragment = \"mappen data som \u00e4r en undermapp under mappen Deluppgift 1.\";\n\n// Locate the paragraph that still contains the outdated folder name and\n// rewrite its full text in one go. Rewriting the whole paragraph (rather\n// than just the matched substring) keeps the resulting run/formatting\n// state simple and avoids leaving stray leftover runs behind.\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (paragraph.text.indexOf(oldFragment) !== -1) {\n    const updatedText = paragraph.text.split(oldFragment).join(newFragment);\n    paragraph.getRange().insertText(updatedText, Word.InsertLocation.replace);\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# The instruction document told students to load data \"fr\u00e5n mappen\n# AdventureWorksDW\" (from the folder AdventureWorksDW). That folder\n# name was wrong; the commit corrects the wording to say the data\n# should be loaded from the \"data\" folder, which is a sub-folder of\n# \"Deluppgift 1\".\n#\n# Old: \"... fr\u00e5n mappen AdventureWorksDW. Transformera data ...\"\n# New: \"... fr\u00e5n mappen data som \u00e4r en undermapp under mappen\n#       Deluppgift 1. Transformera data ...\"\n\n$d = $word.ActiveDocument\n\n$oldFragment = \"mappen AdventureWorksDW.\"\n$newFragment = \"mappen data som \u00e4r en undermapp under mappen Deluppgift 1.\"\n\nforeach ($p in $d.Paragraphs) {\n    $paraText = $p.Range.Text\n    if ($paraText -like \"*$oldFragment*\") {\n        $updatedText = $paraText.Replace($oldFragment, $newFragment)\n        # Rewrite through a document Range (Start/End) rather than the\n        # paragraph's own Range object so the whole paragraph content\n        # (every run) is replaced cleanly in one shot, instead of only\n        # overwriting the first run and leaving the rest behind.\n        $fullRange = $d.Range($p.Range.Start, $p.Range.End)\n        $fullRange.Text = $updatedText\n        break\n    }\n}\n"}
